$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: set the new cell values
$ws.Range("M15").Value = 33.333333333333
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = -35.294117647058
$ws.Range("I16").Value = 106
$ws.Range("J16").Value = 102
$ws.Range("K16").Value = 3.92156862745
$ws.Range("L16").Value = 37.662337662337
$ws.Range("M16").Value = 37.662337662337
$ws.Range("N16").Value = -80.333951762523
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = -33.333333333333
$ws.Range("I17").Value = 93
$ws.Range("J17").Value = 71
$ws.Range("K17").Value = 30.985915492957
$ws.Range("L17").Value = 89.795918367346
$ws.Range("M17").Value = 75.471698113207
$ws.Range("N17").Value = 6.896551724137
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = "'0"
$ws.Range("E18").Value = "'***.*"
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 2
$ws.Range("H18").Value = 150
$ws.Range("I18").Value = 123
$ws.Range("J18").Value = 70
$ws.Range("K18").Value = 75.714285714285
$ws.Range("L18").Value = -8.888888888888
$ws.Range("M18").Value = 25.510204081632
$ws.Range("N18").Value = -84.907975460122
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = 38.461538461538
$ws.Range("F19").Value = 75
$ws.Range("G19").Value = 86
$ws.Range("H19").Value = -12.790697674418
$ws.Range("I19").Value = 752
$ws.Range("J19").Value = 661
$ws.Range("K19").Value = 13.767019667171
$ws.Range("L19").Value = 73.271889400921
$ws.Range("M19").Value = 26.599326599326
$ws.Range("N19").Value = -53.666050523721
$ws.Range("C20").Value = "'0"
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = -50
$ws.Range("I20").Value = 67
$ws.Range("J20").Value = 82
$ws.Range("K20").Value = -18.292682926829
$ws.Range("L20").Value = 55.813953488372
$ws.Range("M20").Value = 131.034482758621
$ws.Range("N20").Value = -93.507751937984
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = 9.090909090909
$ws.Range("F21").Value = 101
$ws.Range("G21").Value = 121
$ws.Range("H21").Value = -16.528925619834
$ws.Range("I21").Value = 1153
$ws.Range("J21").Value = 992
$ws.Range("K21").Value = 16.229838709677
$ws.Range("L21").Value = 54.973118279569
$ws.Range("M21").Value = 33.914053426248
$ws.Range("N21").Value = -71.994170512509
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -75
$ws.Range("I22").Value = 23
$ws.Range("J22").Value = 23
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 27.777777777777
$ws.Range("M22").Value = 4.545454545454
$ws.Range("C23").Value = "'0"
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 0
$ws.Range("L23").Value = 41.176470588235
$ws.Range("C24").Value = 28
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 55.555555555555
$ws.Range("F24").Value = 79
$ws.Range("G24").Value = 78
$ws.Range("H24").Value = 1.282051282051
$ws.Range("I24").Value = 1123
$ws.Range("J24").Value = 1136
$ws.Range("K24").Value = -1.144366197183
$ws.Range("L24").Value = -4.017094017094
$ws.Range("M24").Value = 18.335089567966
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 29
$ws.Range("G25").Value = 27
$ws.Range("H25").Value = 7.407407407407
$ws.Range("I25").Value = 202
$ws.Range("J25").Value = 174
$ws.Range("K25").Value = 16.091954022988
$ws.Range("L25").Value = 66.94214876033
$ws.Range("M25").Value = -19.2
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = "'0"
$ws.Range("E27").Value = "'***.*"
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -60
$ws.Range("I27").Value = 48
$ws.Range("J27").Value = 45
$ws.Range("K27").Value = 6.666666666666
$ws.Range("L27").Value = 50
$ws.Range("F30").Value = "'0"

# Step 2: adjust cell formatting (text vs numeric style) to match the source cells
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4122) | Out-Null
$ws.Range("C14").Copy() | Out-Null
$ws.Range("E18").PasteSpecial(-4122) | Out-Null
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C20").PasteSpecial(-4122) | Out-Null
$ws.Range("F20").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122) | Out-Null
$ws.Range("E19").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4122) | Out-Null
$ws.Range("D22").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4122) | Out-Null
$ws.Range("D27").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$ws.Range("C14").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null
$ws.Range("G27").Copy() | Out-Null
$ws.Range("F27").PasteSpecial(-4122) | Out-Null
$ws.Range("C14").Copy() | Out-Null
$ws.Range("F30").PasteSpecial(-4122) | Out-Null
